# Added R for min load
# Insert a new BOM row (row 18) for a 25-ohm resistor used for the
# minimum-load requirement, pushing the later rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 18 (the blank separator row),
# shifting rows 18-23 down to 19-24.
$ws.Rows("18:18").Insert()

# Populate the new row 18 with the resistor part info.
# (Shared-string insertion order follows the order these are written:
#  Part -> Data Sheet URL -> Part code -> Eagle Code -> Qty.)
$ws.Range("C18").Value = "Resistor (25)"
$ws.Range("G18").Value = "http://www.digikey.com/product-detail/en/504L25R0FTNCFT/1284-1325-1-ND/4789905"
$ws.Range("D18").Value = "504L25R0FTNCFT"
$ws.Range("E18").Value = "R-US_R0402"
$ws.Range("F18").Value = 1

# Move the selection to where the author left it after the edit.
$ws.Range("D19").Select() | Out-Null
